{"js": "// Apply text replacements for the updated date header and multiplication problems.\nconst replacements = [\n  [\"2025-06-07 Saturday\", \"2025-06-08 Sunday\"],\n  [\"71\u00d747=\", \"82\u00d748=\"],\n  [\"97\u00d716=\", \"26\u00d796=\"],\n  [\"79\u00d716=\", \"26\u00d745=\"],\n  [\"42\u00d778=\", \"83\u00d753=\"],\n  [\"94\u00d736=\", \"43\u00d780=\"],\n  [\"29\u00d763=\", \"32\u00d717=\"],\n  [\"88\u00d798=\", \"55\u00d760=\"],\n  [\"19\u00d745=\", \"99\u00d722=\"],\n  [\"87\u00d774=\", \"99\u00d760=\"],\n  [\"67\u00d757=\", \"33\u00d750=\"],\n  [\"88\u00d758=\", \"95\u00d746=\"],\n  [\"62\u00d785=\", \"30\u00d725=\"],\n  [\"95\u00d772=\", \"66\u00d788=\"],\n  [\"92\u00d725=\", \"66\u00d733=\"],\n  [\"88\u00d752=\", \"88\u00d750=\"],\n  [\"52\u00d777=\", \"14\u00d798=\"],\n  [\"64\u00d740=\", \"83\u00d739=\"],\n  [\"44\u00d729=\", \"49\u00d720=\"],\n  [\"87\u00d748=\", \"28\u00d763=\"],\n  [\"70\u00d782=\", \"23\u00d760=\"],\n  [\"64\u00d722=\", \"31\u00d774=\"],\n  [\"51\u00d768=\", \"97\u00d779=\"],\n  [\"56\u00d757=\", \"16\u00d785=\"],\n  [\"19\u00d734=\", \"51\u00d715=\"],\n  [\"89\u00d769=\", \"62\u00d794=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply text replacements for the updated date header and multiplication problems.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"2025-06-07 Saturday\"\n$find.Replacement.Text = \"2025-06-08 Sunday\"\n$find.Execute([ref]\"2025-06-07 Saturday\", $true, $true, $false, $false, $false, $true, 1, $false, \"2025-06-08 Sunday\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"71\u00d747=\"\n$find.Replacement.Text = \"82\u00d748=\"\n$find.Execute([ref]\"71\u00d747=\", $true, $true, $false, $false, $false, $true, 1, $false, \"82\u00d748=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"97\u00d716=\"\n$find.Replacement.Text = \"26\u00d796=\"\n$find.Execute([ref]\"97\u00d716=\", $true, $true, $false, $false, $false, $true, 1, $false, \"26\u00d796=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"79\u00d716=\"\n$find.Replacement.Text = \"26\u00d745=\"\n$find.Execute([ref]\"79\u00d716=\", $true, $true, $false, $false, $false, $true, 1, $false, \"26\u00d745=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"42\u00d778=\"\n$find.Replacement.Text = \"83\u00d753=\"\n$find.Execute([ref]\"42\u00d778=\", $true, $true, $false, $false, $false, $true, 1, $false, \"83\u00d753=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"94\u00d736=\"\n$find.Replacement.Text = \"43\u00d780=\"\n$find.Execute([ref]\"94\u00d736=\", $true, $true, $false, $false, $false, $true, 1, $false, \"43\u00d780=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"29\u00d763=\"\n$find.Replacement.Text = \"32\u00d717=\"\n$find.Execute([ref]\"29\u00d763=\", $true, $true, $false, $false, $false, $true, 1, $false, \"32\u00d717=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"88\u00d798=\"\n$find.Replacement.Text = \"55\u00d760=\"\n$find.Execute([ref]\"88\u00d798=\", $true, $true, $false, $false, $false, $true, 1, $false, \"55\u00d760=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"19\u00d745=\"\n$find.Replacement.Text = \"99\u00d722=\"\n$find.Execute([ref]\"19\u00d745=\", $true, $true, $false, $false, $false, $true, 1, $false, \"99\u00d722=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"87\u00d774=\"\n$find.Replacement.Text = \"99\u00d760=\"\n$find.Execute([ref]\"87\u00d774=\", $true, $true, $false, $false, $false, $true, 1, $false, \"99\u00d760=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"67\u00d757=\"\n$find.Replacement.Text = \"33\u00d750=\"\n$find.Execute([ref]\"67\u00d757=\", $true, $true, $false, $false, $false, $true, 1, $false, \"33\u00d750=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"88\u00d758=\"\n$find.Replacement.Text = \"95\u00d746=\"\n$find.Execute([ref]\"88\u00d758=\", $true, $true, $false, $false, $false, $true, 1, $false, \"95\u00d746=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"62\u00d785=\"\n$find.Replacement.Text = \"30\u00d725=\"\n$find.Execute([ref]\"62\u00d785=\", $true, $true, $false, $false, $false, $true, 1, $false, \"30\u00d725=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"95\u00d772=\"\n$find.Replacement.Text = \"66\u00d788=\"\n$find.Execute([ref]\"95\u00d772=\", $true, $true, $false, $false, $false, $true, 1, $false, \"66\u00d788=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"92\u00d725=\"\n$find.Replacement.Text = \"66\u00d733=\"\n$find.Execute([ref]\"92\u00d725=\", $true, $true, $false, $false, $false, $true, 1, $false, \"66\u00d733=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"88\u00d752=\"\n$find.Replacement.Text = \"88\u00d750=\"\n$find.Execute([ref]\"88\u00d752=\", $true, $true, $false, $false, $false, $true, 1, $false, \"88\u00d750=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"52\u00d777=\"\n$find.Replacement.Text = \"14\u00d798=\"\n$find.Execute([ref]\"52\u00d777=\", $true, $true, $false, $false, $false, $true, 1, $false, \"14\u00d798=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"64\u00d740=\"\n$find.Replacement.Text = \"83\u00d739=\"\n$find.Execute([ref]\"64\u00d740=\", $true, $true, $false, $false, $false, $true, 1, $false, \"83\u00d739=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"44\u00d729=\"\n$find.Replacement.Text = \"49\u00d720=\"\n$find.Execute([ref]\"44\u00d729=\", $true, $true, $false, $false, $false, $true, 1, $false, \"49\u00d720=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"87\u00d748=\"\n$find.Replacement.Text = \"28\u00d763=\"\n$find.Execute([ref]\"87\u00d748=\", $true, $true, $false, $false, $false, $true, 1, $false, \"28\u00d763=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"70\u00d782=\"\n$find.Replacement.Text = \"23\u00d760=\"\n$find.Execute([ref]\"70\u00d782=\", $true, $true, $false, $false, $false, $true, 1, $false, \"23\u00d760=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"64\u00d722=\"\n$find.Replacement.Text = \"31\u00d774=\"\n$find.Execute([ref]\"64\u00d722=\", $true, $true, $false, $false, $false, $true, 1, $false, \"31\u00d774=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"51\u00d768=\"\n$find.Replacement.Text = \"97\u00d779=\"\n$find.Execute([ref]\"51\u00d768=\", $true, $true, $false, $false, $false, $true, 1, $false, \"97\u00d779=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"56\u00d757=\"\n$find.Replacement.Text = \"16\u00d785=\"\n$find.Execute([ref]\"56\u00d757=\", $true, $true, $false, $false, $false, $true, 1, $false, \"16\u00d785=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"19\u00d734=\"\n$find.Replacement.Text = \"51\u00d715=\"\n$find.Execute([ref]\"19\u00d734=\", $true, $true, $false, $false, $false, $true, 1, $false, \"51\u00d715=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"89\u00d769=\"\n$find.Replacement.Text = \"62\u00d794=\"\n$find.Execute([ref]\"89\u00d769=\", $true, $true, $false, $false, $false, $true, 1, $false, \"62\u00d794=\", 2) | Out-Null\n\n"}
